$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize phone numbers in column "tel" (A) that are missing the
# "380" country-code prefix -- replace their first 3 digits with "380".
# (fetchUser / mailing-list data cleanup)
$ws.Range("A6").Value = "380444444444"
$ws.Range("A7").Value = "380555555555"
$ws.Range("A11").Value = "380223333333"
$ws.Range("A12").Value = "380447777777"
$ws.Range("A13").Value = "380555987654"
$ws.Range("A14").Value = "380438459832"
$ws.Range("A18").Value = "380829844444"
$ws.Range("A17").Value = "380333320098"
$ws.Range("A19").Value = "380998055555"
$ws.Range("A23").Value = "380229933333"
$ws.Range("A24").Value = "380777777777"
$ws.Range("A25").Value = "380553437654"
$ws.Range("A5").Value = "380333333333"

# Move the active selection to A5
[void]$ws.Range("A5").Select()
